$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update doctor diagnosis / treatment plan for row 2 (Alice Brown, P1001):
# append ", None" to the existing values
$ws.Range("G2").Value = "require further followup, require medicine, None"
$ws.Range("H2").Value = "followup, done, None"

# Row 3 (Bob Stone, P1002) values stay the same text ("okay, not okay" / "nil, nil")

# Row 4 (Charlie White, P1003) previously had no Doctor Diagnosis / Treatment Plan values,
# now fill them in with "None"
$ws.Range("G4").Value = "None"
$ws.Range("H4").Value = "None"

# Column width tweaks - column F gets an explicit custom width (column G already has
# its custom width set and is left untouched)
$ws.Columns("F").ColumnWidth = 10.5

# Move the active selection to H4
$ws.Range("H4").Select()
